$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "3'5"
$ws.Range("B3").Value = "3'4"
$ws.Range("C4").Value = "KB"
